$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 316.83334
$ws.Cells.Item(55, 9).Value = 180.2
$ws.Cells.Item(55, 11).Value = 180.2
$ws.Cells.Item(55, 13).Value = 33.80000000000001

$ws.Cells.Item(62, 8).Value = 4328.706
$ws.Cells.Item(62, 9).Value = 2184.8572
$ws.Cells.Item(62, 11).Value = 2184.8572
$ws.Cells.Item(62, 13).Value = -1560.8572

$ws.Cells.Item(65, 8).Value = 4328.706
$ws.Cells.Item(65, 9).Value = 2184.8572
$ws.Cells.Item(65, 11).Value = 10924.286
$ws.Cells.Item(65, 13).Value = -7804.286

$ws.Cells.Item(70, 8).Value = 4633.3335
$ws.Cells.Item(70, 9).Value = 9800
$ws.Cells.Item(70, 10).Value = 2050
$ws.Cells.Item(70, 11).Value = 29400
$ws.Cells.Item(70, 12).Value = 6150
$ws.Cells.Item(70, 13).Value = -29130
$ws.Cells.Item(70, 14).Value = -6690

$ws.Cells.Item(73, 8).Value = 4633.3335
$ws.Cells.Item(73, 9).Value = 9800
$ws.Cells.Item(73, 10).Value = 2050
$ws.Cells.Item(73, 11).Value = 29400
$ws.Cells.Item(73, 12).Value = 6150
$ws.Cells.Item(73, 13).Value = -28464
$ws.Cells.Item(73, 14).Value = -8022

$ws.Cells.Item(88, 8).Value = 4024.7368
$ws.Cells.Item(88, 9).Value = 6450
$ws.Cells.Item(88, 10).Value = 2260.9092
$ws.Cells.Item(88, 11).Value = 6450
$ws.Cells.Item(88, 12).Value = 2260.9092
$ws.Cells.Item(88, 13).Value = -6044
$ws.Cells.Item(88, 14).Value = -3072.9092

$ws.Cells.Item(91, 8).Value = 4024.7368
$ws.Cells.Item(91, 9).Value = 6450
$ws.Cells.Item(91, 10).Value = 2260.9092
$ws.Cells.Item(91, 11).Value = 6450
$ws.Cells.Item(91, 12).Value = 2260.9092
$ws.Cells.Item(91, 13).Value = -5046
$ws.Cells.Item(91, 14).Value = -5068.9092

$ws.Cells.Item(137, 8).Value = 1737.4572
$ws.Cells.Item(137, 9).Value = 1130.5652
$ws.Cells.Item(137, 10).Value = 2900.6667
$ws.Cells.Item(137, 11).Value = 3391.6956
$ws.Cells.Item(137, 12).Value = 8702.000100000001
$ws.Cells.Item(137, 13).Value = -841.6956
$ws.Cells.Item(137, 14).Value = -13802.0001

$ws.Cells.Item(138, 8).Value = 118171.03
$ws.Cells.Item(138, 9).Value = 244493.12
$ws.Cells.Item(138, 10).Value = 3077.578
$ws.Cells.Item(138, 11).Value = 733479.36
$ws.Cells.Item(138, 12).Value = 9232.734
$ws.Cells.Item(138, 13).Value = -728339.36
$ws.Cells.Item(138, 14).Value = -19512.734

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1149.8823
$ws.Cells.Item(61, 9).Value = 1034.475
$ws.Cells.Item(61, 10).Value = 1569.5454
$ws.Cells.Item(61, 11).Value = 1034.475
$ws.Cells.Item(61, 12).Value = 1569.5454
$ws.Cells.Item(61, 13).Value = -822.4749999999999
$ws.Cells.Item(61, 14).Value = -1993.5454

$ws.Cells.Item(88, 8).Value = 500001500
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 500001500
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 500001500
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(88, 14).Value = -500002312

$ws.Cells.Item(91, 8).Value = 500001500
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 500001500
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 500001500
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(91, 14).Value = -500004308

$ws.Cells.Item(136, 8).Value = 1149.8823
$ws.Cells.Item(136, 9).Value = 1034.475
$ws.Cells.Item(136, 10).Value = 1569.5454
$ws.Cells.Item(136, 11).Value = 3103.425
$ws.Cells.Item(136, 12).Value = 4708.6362
$ws.Cells.Item(136, 13).Value = -553.4249999999997
$ws.Cells.Item(136, 14).Value = -9808.636200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2813.9333
$ws.Cells.Item(20, 9).Value = 3377
$ws.Cells.Item(20, 10).Value = 2609.182
$ws.Cells.Item(20, 11).Value = 3377
$ws.Cells.Item(20, 12).Value = 2609.182
$ws.Cells.Item(20, 13).Value = -3130
$ws.Cells.Item(20, 14).Value = -3103.182

$ws.Cells.Item(80, 8).Value = 43.8125
$ws.Cells.Item(80, 10).Value = 49.214287
$ws.Cells.Item(80, 12).Value = 49.214287
$ws.Cells.Item(80, 14).Value = -2045.214287

$ws.Cells.Item(83, 8).Value = 43.8125
$ws.Cells.Item(83, 10).Value = 49.214287
$ws.Cells.Item(83, 12).Value = 246.071435
$ws.Cells.Item(83, 14).Value = -10230.071435

$ws.Cells.Item(86, 8).Value = 13335420
$ws.Cells.Item(86, 9).Value = 25001776
$ws.Cells.Item(86, 10).Value = 2442.8572
$ws.Cells.Item(86, 11).Value = 25001776
$ws.Cells.Item(86, 12).Value = 2442.8572
$ws.Cells.Item(86, 13).Value = -25000653
$ws.Cells.Item(86, 14).Value = -4688.8572

$ws.Cells.Item(89, 8).Value = 13335420
$ws.Cells.Item(89, 9).Value = 25001776
$ws.Cells.Item(89, 10).Value = 2442.8572
$ws.Cells.Item(89, 11).Value = 125008880
$ws.Cells.Item(89, 12).Value = 12214.286
$ws.Cells.Item(89, 13).Value = -125003264
$ws.Cells.Item(89, 14).Value = -23446.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1126.6415
$ws.Cells.Item(58, 9).Value = 719.3913
$ws.Cells.Item(58, 10).Value = 1438.8667
$ws.Cells.Item(58, 11).Value = 719.3913
$ws.Cells.Item(58, 12).Value = 1438.8667
$ws.Cells.Item(58, 13).Value = -516.3913
$ws.Cells.Item(58, 14).Value = -1844.8667

$ws.Cells.Item(94, 8).Value = 250000850
$ws.Cells.Item(94, 9).Value = 500000540
$ws.Cells.Item(94, 11).Value = 500000540
$ws.Cells.Item(94, 13).Value = -500000089

$ws.Cells.Item(99, 8).Value = 1154599.4
$ws.Cells.Item(99, 9).Value = 2103575
$ws.Cells.Item(99, 11).Value = 2103575
$ws.Cells.Item(99, 13).Value = -2102077

$ws.Cells.Item(126, 8).Value = 1154599.4
$ws.Cells.Item(126, 9).Value = 2103575
$ws.Cells.Item(126, 11).Value = 6310725
$ws.Cells.Item(126, 13).Value = -6308255

$ws.Cells.Item(136, 8).Value = 1126.6415
$ws.Cells.Item(136, 9).Value = 719.3913
$ws.Cells.Item(136, 10).Value = 1438.8667
$ws.Cells.Item(136, 11).Value = 2158.1739
$ws.Cells.Item(136, 12).Value = 4316.6001
$ws.Cells.Item(136, 13).Value = 391.8261000000002
$ws.Cells.Item(136, 14).Value = -9416.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 584.8946999999999
$ws.Cells.Item(5, 9).Value = 547.5
$ws.Cells.Item(5, 10).Value = 689.6
$ws.Cells.Item(5, 11).Value = 1642.5
$ws.Cells.Item(5, 12).Value = 2068.8
$ws.Cells.Item(5, 13).Value = -1530.5
$ws.Cells.Item(5, 14).Value = -2292.8

$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 14).ClearContents()

$ws.Cells.Item(109, 8).Value = 4063.158
$ws.Cells.Item(109, 9).Value = 2333.1667
$ws.Cells.Item(109, 10).Value = 4861.615
$ws.Cells.Item(109, 11).Value = 6999.500100000001
$ws.Cells.Item(109, 12).Value = 14584.845
$ws.Cells.Item(109, 13).Value = -5959.500100000001
$ws.Cells.Item(109, 14).Value = -16664.845

$ws.Cells.Item(135, 8).Value = 584.8946999999999
$ws.Cells.Item(135, 9).Value = 547.5
$ws.Cells.Item(135, 10).Value = 689.6
$ws.Cells.Item(135, 11).Value = 4927.5
$ws.Cells.Item(135, 12).Value = 6206.400000000001
$ws.Cells.Item(135, 13).Value = -2392.5
$ws.Cells.Item(135, 14).Value = -11276.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 57147772
$ws.Cells.Item(70, 9).Value = 133337770
$ws.Cells.Item(70, 11).Value = 133337770
$ws.Cells.Item(70, 13).Value = -133337500

$ws.Cells.Item(73, 8).Value = 57147772
$ws.Cells.Item(73, 9).Value = 133337770
$ws.Cells.Item(73, 11).Value = 133337770
$ws.Cells.Item(73, 13).Value = -133336834

$ws.Cells.Item(80, 8).Value = 3532.1785
$ws.Cells.Item(80, 9).Value = 6541
$ws.Cells.Item(80, 10).Value = 2878.087
$ws.Cells.Item(80, 11).Value = 6541
$ws.Cells.Item(80, 12).Value = 2878.087
$ws.Cells.Item(80, 13).Value = -5543
$ws.Cells.Item(80, 14).Value = -4874.087

$ws.Cells.Item(83, 8).Value = 3532.1785
$ws.Cells.Item(83, 9).Value = 6541
$ws.Cells.Item(83, 10).Value = 2878.087
$ws.Cells.Item(83, 11).Value = 32705
$ws.Cells.Item(83, 12).Value = 14390.435
$ws.Cells.Item(83, 13).Value = -27713
$ws.Cells.Item(83, 14).Value = -24374.435

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 329.18182
$ws.Cells.Item(55, 9).Value = 323.5
$ws.Cells.Item(55, 10).Value = 336
$ws.Cells.Item(55, 11).Value = 323.5
$ws.Cells.Item(55, 12).Value = 336
$ws.Cells.Item(55, 13).Value = -150.5
$ws.Cells.Item(55, 14).Value = -682

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 90912024
$ws.Cells.Item(81, 9).Value = 166669570
$ws.Cells.Item(81, 10).Value = 2980
$ws.Cells.Item(81, 11).Value = 333339140
$ws.Cells.Item(81, 12).Value = 5960
$ws.Cells.Item(81, 13).Value = -333338079
$ws.Cells.Item(81, 14).Value = -8082

$ws.Cells.Item(84, 8).Value = 90912024
$ws.Cells.Item(84, 9).Value = 166669570
$ws.Cells.Item(84, 10).Value = 2980
$ws.Cells.Item(84, 11).Value = 1666695700
$ws.Cells.Item(84, 12).Value = 29800
$ws.Cells.Item(84, 13).Value = -1666690396
$ws.Cells.Item(84, 14).Value = -40408

$ws.Cells.Item(126, 8).Value = 5800.5
$ws.Cells.Item(126, 9).Value = 8580.799999999999
$ws.Cells.Item(126, 11).Value = 25742.4
$ws.Cells.Item(126, 13).Value = -23272.4
